$wb = $excel.ActiveWorkbook

# --- Debts sheet: insert a new "active" column at column A ---
$debts = $wb.Worksheets.Item("Debts")
$debts.Activate()
$debts.Range("A1").EntireColumn.Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)
$debts.Range("A1").Value = "active"
$debts.Range("A1").Font.Bold = $true

# --- Fixed Assets sheet: insert a new "active" column at column A ---
$assets = $wb.Worksheets.Item("Fixed Assets")
$assets.Activate()
$assets.Range("A1").EntireColumn.Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)
$assets.Range("A1").Value = "active"
$assets.Range("A1").Font.Bold = $true

# Final selection / active sheet state matches the authored edit: Fixed Assets
# becomes the active tab, with C11 selected, while Debts is left with column A
# selected (the column that was just inserted).
$debts.Range("A1:A1048576").Select()
$assets.Activate()
$assets.Range("C11").Select()
